# Generate Report for Handoff
#
# The source file 813d25a4-d47d-4158-9778-15ee23e06a8b.md has now been
# handed off for translation (zh-cn, de-de): status moves from
# "In Translation" to "Ready for handoff", priority from "ht" to "mt",
# and the handoff timestamps / overview generation date are refreshed.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for 813d25a4-...md ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-19 02:13:43"
$wsOverview.Columns.Item(5).ColumnWidth = 16.3
$wsOverview.Columns.Item(6).ColumnWidth = 16.3

# --- zh-cn sheet: row for 813d25a4-...md ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-19 02:13:38"
$wsZhCn.Columns.Item(3).ColumnWidth = 16.3

# --- de-de sheet: row for 813d25a4-...md ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-19 02:13:43"
$wsDeDe.Columns.Item(3).ColumnWidth = 16.3
